$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the month/date column (B) for all data rows (2-20) to the new period (1-Apr-2023, serial 45017)
$ws.Range("B2:B20").Value = 45017

# Update Cantidad (D) values that changed for the new month's data
$ws.Cells.Item(2, 4).Value = 99
$ws.Cells.Item(3, 4).Value = 134
$ws.Cells.Item(4, 4).Value = 19
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(7, 4).Value = 7
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(9, 4).Value = 1

# Row 10 (San Martin) procedure changed from CURACION to GLUCOMETRIA, with a new quantity
$ws.Cells.Item(10, 3).Value = "GLUCOMETRIA"
$ws.Cells.Item(10, 4).Value = 60

$ws.Cells.Item(13, 4).Value = 123
$ws.Cells.Item(15, 4).Value = 57
$ws.Cells.Item(16, 4).Value = 56
$ws.Cells.Item(17, 4).Value = 57

# The old row 20 (Cartagena / CURACION) no longer exists in the new month's data - remove it
$ws.Rows.Item(20).Delete()

# Move the selection to reflect where the user left off editing
$ws.Range("D14").Select()
